$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Census")

# Clear the stale note that does not carry over once row 16 is reshaped
$ws.Range("G16").ClearContents()

$ws.Range("A1").Value = "Outcome Category"
$ws.Range("B1").Value = "Outcome"
$ws.Range("C1").Value = "Dataset"
$ws.Range("D1").Value = "Main Dataset"
$ws.Range("E1").Value = "Geography"
$ws.Range("F1").Value = "Years"
$ws.Range("G1").Value = "Notes"
$ws.Range("H1").Value = "Downloaded"

$ws.Range("A2").Value = "Jobs creation"
$ws.Range("B2").Value = "# establishments"
$ws.Range("C2").Value = "Zip Code Business Patterns (ZBP)"
$ws.Range("D2").Value = "Zip Code Business Patterns (ZBP)"
$ws.Range("E2").Value = "zip code"
$ws.Range("F2").Value = "1994-2016"
$ws.Range("G2").Value = "week of March 12, also gives Q1 and annual payroll. Gives it by SIC/NAICS, but we'll just download aggregate for zip code."
$ws.Range("H2").Value = "x"

$ws.Range("A3").Value = "Jobs creation"
$ws.Range("B3").Value = "# employed"
$ws.Range("C3").Value = "Zip Code Business Patterns (ZBP)"
$ws.Range("D3").Value = "Zip Code Business Patterns (ZBP)"
$ws.Range("E3").Value = "zip code"
$ws.Range("F3").Value = "1994-2016"
$ws.Range("G3").Value = "week of March 12, also gives Q1 and annual payroll. Gives it by SIC/NAICS, but we'll just download aggregate for zip code."
$ws.Range("H3").Value = "x"

$ws.Range("A4").Value = "Jobs creation"
$ws.Range("B4").Value = "annual payroll"
$ws.Range("C4").Value = "Zip Code Business Patterns (ZBP)"
$ws.Range("D4").Value = "Zip Code Business Patterns (ZBP)"
$ws.Range("E4").Value = "zip code"
$ws.Range("F4").Value = "1994-2016"
$ws.Range("G4").Value = "week of March 12, also gives Q1 and annual payroll. Gives it by SIC/NAICS, but we'll just download aggregate for zip code."
$ws.Range("H4").Value = "x"

$ws.Range("A5").Value = "Median wages"
$ws.Range("B5").Value = "median income by race/ethnicity"
$ws.Range("C5").Value = "ACS Median Income in Past 12 months (S1903)"
$ws.Range("D5").Value = "ACS 5 yr Subject"
$ws.Range("E5").Value = "tract"
$ws.Range("F5").Value = "2010-2017"
$ws.Range("H5").Value = "x"

$ws.Range("A6").Value = "Median wages"
$ws.Range("B6").Value = "median income by age group"
$ws.Range("C6").Value = "ACS Median Income in Past 12 months (S1903)"
$ws.Range("D6").Value = "ACS 5 yr Subject"
$ws.Range("E6").Value = "tract"
$ws.Range("F6").Value = "2010-2017"
$ws.Range("H6").Value = "??"

$ws.Range("A7").Value = "Education"
$ws.Range("B7").Value = "school enrollment (# kids enrolled by school level)"
$ws.Range("C7").Value = "ACS School Enrollment (S1401)"
$ws.Range("D7").Value = "ACS 5 yr Subject"
$ws.Range("E7").Value = "tract"
$ws.Range("F7").Value = "2010-2017"
$ws.Range("H7").Value = "??"

$ws.Range("A8").Value = "Education"
$ws.Range("B8").Value = "educational attainment by age groups"
$ws.Range("C8").Value = "ACS Educational Attainment (S1501)"
$ws.Range("D8").Value = "ACS 5 yr Subject"
$ws.Range("E8").Value = "tract"
$ws.Range("F8").Value = "2010-2017"
$ws.Range("G8").Value = "only for pop 25+, did not include 18-24"
$ws.Range("H8").Value = "x"

$ws.Range("A9").Value = "Education"
$ws.Range("B9").Value = "educational attainment by race"
$ws.Range("C9").Value = "ACS Educational Attainment (S1501)"
$ws.Range("D9").Value = "ACS 5 yr Subject"
$ws.Range("E9").Value = "tract"
$ws.Range("F9").Value = "2010-2017"
$ws.Range("H9").Value = "x"

$ws.Range("A10").Value = "Quality of life"
$ws.Range("B10").Value = "poverty rate for 25+ by educational attainment"
$ws.Range("C10").Value = "ACS Educational Attainment (S1501)"
$ws.Range("D10").Value = "ACS 5 yr Subject"
$ws.Range("E10").Value = "tract"
$ws.Range("F10").Value = "2010-2017"
$ws.Range("H10").Value = "x"

$ws.Range("A11").Value = "Median wages"
$ws.Range("B11").Value = "median earnings by educational attainment"
$ws.Range("C11").Value = "ACS Educational Attainment (S1501)"
$ws.Range("D11").Value = "ACS 5 yr Subject"
$ws.Range("E11").Value = "tract"
$ws.Range("F11").Value = "2010-2017"
$ws.Range("H11").Value = "x"

$ws.Range("A12").Value = "Quality of life"
$ws.Range("B12").Value = "poverty status by race"
$ws.Range("C12").Value = "ACS Poverty Status in past 12 Months (S1701)"
$ws.Range("D12").Value = "ACS 5 yr Subject"
$ws.Range("E12").Value = "tract"
$ws.Range("F12").Value = "2012-2017"
$ws.Range("H12").Value = "x"

$ws.Range("A13").Value = "Quality of life"
$ws.Range("B13").Value = "poverty status by educational attainment"
$ws.Range("C13").Value = "ACS Poverty Status in past 12 Months (S1701)"
$ws.Range("D13").Value = "ACS 5 yr Subject"
$ws.Range("E13").Value = "tract"
$ws.Range("F13").Value = "2012-2017"
$ws.Range("H13").Value = "x"

$ws.Range("A14").Value = "Quality of life"
$ws.Range("B14").Value = "poverty status of families"
$ws.Range("C14").Value = "ACS Poverty Status in past 12 Months (S1702)"
$ws.Range("D14").Value = "ACS 5 yr Subject"
$ws.Range("E14").Value = "tract"
$ws.Range("F14").Value = "2010-2017"
$ws.Range("H14").Value = "x"

$ws.Range("A15").Value = "Quality of life"
$ws.Range("B15").Value = "households receiving food stamps/SNAP"
$ws.Range("C15").Value = "ACS Food Stamps/SNAP (S2201)"
$ws.Range("D15").Value = "ACS 5 yr Subject"
$ws.Range("E15").Value = "tract"
$ws.Range("F15").Value = "2010-2017"
$ws.Range("H15").Value = "x"

$ws.Range("A16").Value = "Quality of life"
$ws.Range("B16").Value = "households receiving food stamps/SNAP by poverty level"
$ws.Range("C16").Value = "ACS Food Stamps/SNAP (S2201)"
$ws.Range("D16").Value = "ACS 5 yr Subject"
$ws.Range("E16").Value = "tract"
$ws.Range("F16").Value = "2010-2017"
$ws.Range("H16").Value = "x"

$ws.Range("A17").Value = "Employment rates"
$ws.Range("B17").Value = "employment status by age group"
$ws.Range("C17").Value = "ACS Employment Status (S2301)"
$ws.Range("D17").Value = "ACS 5 yr Subject"
$ws.Range("E17").Value = "tract"
$ws.Range("F17").Value = "2010-2017"
$ws.Range("G17").Value = "employment status measured by labor force participation rate, employment/population ratio, and unemployment rate"
$ws.Range("H17").Value = "??"

$ws.Range("A18").Value = "Employment rates"
$ws.Range("B18").Value = "employment status by race/ethnicity"
$ws.Range("C18").Value = "ACS Employment Status (S2301)"
$ws.Range("D18").Value = "ACS 5 yr Subject"
$ws.Range("E18").Value = "tract"
$ws.Range("F18").Value = "2010-2017"
$ws.Range("G18").Value = "employment status measured by labor force participation rate, employment/population ratio, and unemployment rate"
$ws.Range("H18").Value = "x"

$ws.Range("A19").Value = "Employment rates"
$ws.Range("B19").Value = "employment status by sex"
$ws.Range("C19").Value = "ACS Employment Status (S2301)"
$ws.Range("D19").Value = "ACS 5 yr Subject"
$ws.Range("E19").Value = "tract"
$ws.Range("F19").Value = "2010-2017"
$ws.Range("G19").Value = "employment status measured by labor force participation rate, employment/population ratio, and unemployment rate"
$ws.Range("H19").Value = "x"

$ws.Range("A20").Value = "Employment rates"
$ws.Range("B20").Value = "employment status by poverty level"
$ws.Range("C20").Value = "ACS Employment Status (S2301)"
$ws.Range("D20").Value = "ACS 5 yr Subject"
$ws.Range("E20").Value = "tract"
$ws.Range("F20").Value = "2010-2017"
$ws.Range("G20").Value = "employment status measured by labor force participation rate, employment/population ratio, and unemployment rate"
$ws.Range("H20").Value = "x"

$ws.Range("A21").Value = "Employment rates"
$ws.Range("B21").Value = "employment status by educational attainment"
$ws.Range("C21").Value = "ACS Employment Status (S2301)"
$ws.Range("D21").Value = "ACS 5 yr Subject"
$ws.Range("E21").Value = "tract"
$ws.Range("F21").Value = "2010-2017"
$ws.Range("G21").Value = "employment status measured by labor force participation rate, employment/population ratio, and unemployment rate"
$ws.Range("H21").Value = "x"

$ws.Range("A22").Value = "Quality of life"
$ws.Range("B22").Value = "monthly housing cost as % of household income"
$ws.Range("C22").Value = "ACS Financial Characteristics (S2503)"
$ws.Range("D22").Value = "ACS 5 yr Subject"
$ws.Range("E22").Value = "tract"
$ws.Range("F22").Value = "2010-2017"
$ws.Range("G22").Value = "categorical: <20, 20-30, 30+"
$ws.Range("H22").Value = "??"

$ws.Range("A23").Value = "Quality of life"
$ws.Range("B23").Value = "% of hh that receive public assistance"
$ws.Range("C23").Value = "ACS Receipt of Supplemental Security Income (SSI), Cash Public Assistance Income, or Food Stamps/SNAPS for children under 18 (B09010)"
$ws.Range("D23").Value = "ACS 5 yr"
$ws.Range("E23").Value = "tract"
$ws.Range("F23").Value = "2010-2017"
$ws.Range("G23").Value = "universe: population under 18 yrs in hh"
$ws.Range("H23").Value = "??"

$ws.Range("A24").Value = "Quality of life"
$ws.Range("B24").Value = "% of households that receive public assistance income or food stamps/SNAP"
$ws.Range("C24").Value = "ACS Public Assistance Income or Food Stamps/SNAP in past 12 Months (B19058)"
$ws.Range("D24").Value = "ACS 5 yr"
$ws.Range("E24").Value = "tract"
$ws.Range("F24").Value = "2010-2017"
$ws.Range("G24").Value = "universe: households"
$ws.Range("H24").Value = "x"

$ws.Range("A25").Value = "Quality of life"
$ws.Range("B25").Value = "aggregate public assistance income for past 12 months"
$ws.Range("C25").Value = "ACS Aggregate Public Assistance Income in Past 12 Months (B19067)"
$ws.Range("D25").Value = "ACS 5 yr"
$ws.Range("E25").Value = "tract"
$ws.Range("F25").Value = "2010-2017"
$ws.Range("G25").Value = "will need to normalize by population"
$ws.Range("H25").Value = "x"

$ws.Range("A26").Value = "Population"
$ws.Range("B26").Value = "total population"
$ws.Range("C26").Value = "ACS Total Population (B01003)"
$ws.Range("D26").Value = "ACS 5 yr"
$ws.Range("E26").Value = "tract"
$ws.Range("F26").Value = "2010-2017"
$ws.Range("G26").Value = "universe: total population. not explicitly needed, but use to normalize any aggregate values"
$ws.Range("H26").Value = "x"

$ws.Range("A27").Value = "Housing"
$ws.Range("B27").Value = "total housing units"
$ws.Range("C27").Value = "ACS Housing Units (B25001)"
$ws.Range("D27").Value = "ACS 5 yr"
$ws.Range("E27").Value = "tract"
$ws.Range("F27").Value = "2010-2017"
$ws.Range("G27").Value = "universe: housing units. not explicitly needed, but use to normalize any aggregate values"
$ws.Range("H27").Value = "x"

$ws.Range("A28").Value = "Population"
$ws.Range("B28").Value = "total population"
$ws.Range("C28").Value = "ACS Total Population (B01003)"
$ws.Range("D28").Value = "ACS 5 yr"
$ws.Range("E28").Value = "block groups"
$ws.Range("F28").Value = "2013-2017"
$ws.Range("G28").Value = "not explicitly needed, but might use this to aggregate tracts to congressional districts or zip codes, especially in how to split census tracts if they don't fall neatly "

# Update the frozen-pane anchor and active selection to match the new layout
$ws.Activate()
$ws.Range("G7").Select()
